# Insert a new "FindCarTest" worksheet between "LoginTest" and "UserRegTest",
# and populate it with car-brand / browser-type / run-mode test data.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() inserts the new sheet right before the (currently active)
# sheet - since UserRegTest (index 1, "activeTab=1") is active, this lands
# the new sheet between LoginTest and UserRegTest, exactly like the diff.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "FindCarTest"

# Header row.
$newSheet.Range("A1").Value = "carbrand"
$newSheet.Range("B1").Value = "browserType"
$newSheet.Range("C1").Value = "runmode"

# Fill column-by-column (matches shared-string insertion order seen upstream:
# headers, then column B's distinct values, then column C's, then column A's).
$newSheet.Range("B2").Value = "chrome"
$newSheet.Range("B3").Value = "firefox"
$newSheet.Range("B4").Value = "chrome"

$newSheet.Range("C2").Value = "Y"
$newSheet.Range("C3").Value = "N"
$newSheet.Range("C4").Value = "Y"

$newSheet.Range("A2").Value = "bmw"
$newSheet.Range("A3").Value = "mg"
$newSheet.Range("A4").Value = "toyota"

# Widen column B slightly (stored width ends up at 14 once Excel's internal
# character-width padding is applied).
$newSheet.Columns.Item(2).ColumnWidth = (158/12)

# Leave the selection one row below the data, as a human entering this table
# and pressing Enter down column A would.
$newSheet.Range("A5").Select() | Out-Null
